$wb = $excel.ActiveWorkbook

# --- Users sheet: shift existing user into column B, add new user in A2 ---
$usersWs = $wb.Worksheets.Item("Users")
$existingUser = $usersWs.Range("A2").Value2
$usersWs.Range("B2").Value = $existingUser
$usersWs.Range("B2").WrapText = $true
$usersWs.Columns.Item(2).ColumnWidth = 12.6
$usersWs.Range("A2").Value = "Julie Carthane"

# --- GiftLog sheet: update submitted-for user to the new name ---
$giftLogWs = $wb.Worksheets.Item("GiftLog")
$giftLogWs.Range("B2").Value = "Julie Carthane"

# move selection on GiftLog sheet
$giftLogWs.Activate()
$giftLogWs.Range("H7").Select()
